# Update the Nxph3-Nrxn2 LR-pairs sheet with the new TPM-derived values.
# The underlying data now also includes a new "Inflammatory-Mac" target
# cluster for both sending clusters (ECs, MuSCs), growing the table from
# 6 data rows to 8 data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..T the numeric metric columns (same order as the header row).
$data = @(
  @("ECs",   "Nxph3", "Nrxn2", "ECs",              1, 0.3333333333333333, 0.06397866666666667,   0.191936, 0.8900471604056629, 0.8900471604056629, 3, 1,                  0.173461,            0.520383, 0.2181774959583017,  0.2181774959583018,  0.01109780349866667,   0.099880231488,    0.1941882607421045,  0.1941882607421045),
  @("ECs",   "Nxph3", "Nrxn2", "FAPs",             1, 0.3333333333333333, 0.06397866666666667,   0.191936, 0.8900471604056629, 0.8900471604056629, 3, 1,                  0.4903776666666667, 1.471133, 0.6167920822963554,  0.6167920822963555,  0.03137370927644444,   0.282363383488,    0.5489740414085671,  0.5489740414085672),
  @("ECs",   "Nxph3", "Nrxn2", "Inflammatory-Mac", 1, 0.3333333333333333, 0.06397866666666667,   0.191936, 0.8900471604056629, 0.8900471604056629, 1, 0.3333333333333333, 0.04417666666666667, 0.13253,  0.0555649656874912,  0.0555649656874912,  0.002826364231111112,  0.02543727808,     0.04945543992818963, 0.04945543992818964),
  @("ECs",   "Nxph3", "Nrxn2", "MuSCs",            1, 0.3333333333333333, 0.06397866666666667,   0.191936, 0.8900471604056629, 0.8900471604056629, 3, 1,                  0.08703,             0.26109,  0.1094654560578516,  0.1094654560578516,  0.00556806336,         0.05011257023999999, 0.0974294183268017, 0.09742941832680173),
  @("MuSCs", "Nxph3", "Nrxn2", "ECs",              1, 0.3333333333333333, 0.007903666666666666,  0.023711, 0.109952839594337,  0.109952839594337,  3, 1,                  0.173461,            0.520383, 0.2181774959583017,  0.2181774959583018,  0.001370977923666667,  0.012338801313,    0.02398923521619727, 0.02398923521619727),
  @("MuSCs", "Nxph3", "Nrxn2", "FAPs",             1, 0.3333333333333333, 0.007903666666666666,  0.023711, 0.109952839594337,  0.109952839594337,  3, 1,                  0.4903776666666667, 1.471133, 0.6167920822963554,  0.6167920822963555,  0.003875781618111111,  0.034882034563,    0.0678180408877883,  0.06781804088778831),
  @("MuSCs", "Nxph3", "Nrxn2", "Inflammatory-Mac", 1, 0.3333333333333333, 0.007903666666666666,  0.023711, 0.109952839594337,  0.109952839594337,  1, 0.3333333333333333, 0.04417666666666667, 0.13253,  0.0555649656874912,  0.0555649656874912,  0.0003491576477777778, 0.00314241883,     0.006109525759301561, 0.006109525759301562),
  @("MuSCs", "Nxph3", "Nrxn2", "MuSCs",            1, 0.3333333333333333, 0.007903666666666666,  0.023711, 0.109952839594337,  0.109952839594337,  3, 1,                  0.08703,             0.26109,  0.1094654560578516,  0.1094654560578516,  0.0006878561099999999, 0.00619070499,     0.01203603773104991, 0.01203603773104991)
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $rowNum = $i + 2
  $rowVals = $data[$i]
  for ($j = 0; $j -lt $rowVals.Length; $j++) {
    $colNum = $j + 1
    $ws.Cells.Item($rowNum, $colNum).Value = $rowVals[$j]
  }
}
